# Updated symbol list on Sun Dec 18 21:31:07 UTC 2022 with GitHub Actions
#
# This script applies refreshed coin price / volume-label values to the
# "cryptos" worksheet, mirroring a scheduled scraper update.
#
# Price (column D) cells hold numeric-looking text (they are stored as
# strings in the workbook, not numbers), so a leading apostrophe is used
# to force Excel to keep them as text instead of silently converting them
# to numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
}

# Column D (Price) updates
Set-TextValue "D2"  "248.29"
Set-TextValue "D4"  "5.562"
Set-TextValue "D5"  "0.05667"
Set-TextValue "D6"  "6.443"
Set-TextValue "D7"  "0.8006"
Set-TextValue "D9"  "0.1434"
Set-TextValue "D10" "0.07327"
Set-TextValue "D11" "0.03130"
Set-TextValue "D12" "0.02922"
Set-TextValue "D13" "0.09274"
Set-TextValue "D14" "0.001658"
Set-TextValue "D15" "3.220"
Set-TextValue "D16" "0.04745"
Set-TextValue "D17" "0.0005816"
Set-TextValue "D19" "0.005062"
Set-TextValue "D20" "0.001050"
Set-TextValue "D21" "0.0001501"
Set-TextValue "D22" "3.976"
Set-TextValue "D24" "2.088"
Set-TextValue "D25" "0.3267"
Set-TextValue "D26" "0.1265"
Set-TextValue "D27" "0.0003203"
Set-TextValue "D40" "0.04128"
Set-TextValue "D41" "0.007007"
Set-TextValue "D42" "0.003503"
Set-TextValue "D43" "0.1042"
Set-TextValue "D44" "0.009402"
Set-TextValue "D45" "0.00005643"
Set-TextValue "D47" "0.7860"
Set-TextValue "D48" "0.01646"
Set-TextValue "D49" "0.00002102"

# Column E (Volume(1h) label) updates - plain text, no quoting needed
$ws.Range("E17").Value = "16OneONEWorstin24h"
$ws.Range("E42").Value = "41CEJICEJIBestin24h"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("E48").Value = "47BOLOBOLO"
